$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (the 2nd paragraph, right after the title).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Replace the text of the final "Prompt: ..." paragraph with the meta-description text,
#    keeping its existing (italic) run formatting.
$found = $d.Content.Find.Execute(
    "Prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses to fit the theme of the online slot game " + [char]34 + "1 Left Alive." + [char]34 + " The image should be colorful and depict the Maya warrior holding a weapon and surrounded by zombies. The warrior should be wearing traditional Maya clothing and a headpiece, and the glasses should be prominent. Make sure that the image is eye-catching and will draw in players who enjoy action-packed slot games.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the top features of 1 Left Alive slot game, including high payout percentage and special symbols. Play for free or real money.",
    2
)

# 3. Insert a new bold paragraph, "Play 1 Left Alive Slot Game for Free - Review and Features",
#    right before that last paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$null = $lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 1 Left Alive Slot Game for Free - Review and Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newPara.Range.InsertXML($xml)
